# Bump the "Förändrad" (Changed) date in column C by one day for all
# data rows (rows 2 through 236) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 236

$range = $ws.Range("C$firstRow`:C$lastRow")
$range.Value = 45190
